$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.745.27'
$ws.Range("E2").Value = '  +2.27%  '
$ws.Range("D3").Value = '1.874.83'
$ws.Range("E3").Value = '  +2.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.08'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4597'
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3862'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9945'
$ws.Range("E10").Value = '  +3.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.73'
$ws.Range("E11").Value = '  -1.13%  '
$ws.Range("D12").Value = '1.892.76'
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.984'
$ws.Range("E13").Value = '  +1.30%  '
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06964'
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.44'
$ws.Range("E16").Value = '  +0.22%  '
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001006'
$ws.Range("E18").Value = '  +1.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.84'
$ws.Range("E19").Value = '  +0.83%  '
$ws.Range("E20").Value = '  +0.25%  '
$ws.Range("D21").Value = '28.745.92'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.126'
$ws.Range("E24").Value = '  +1.78%  '
$ws.Range("D25").Value = '2.124.26'
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.46'
$ws.Range("E26").Value = '  -0.77%  '
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.777'
$ws.Range("E28").Value = '  +0.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.959'
$ws.Range("E29").Value = '  -0.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.90'
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("E31").Value = '  +0.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9185'
$ws.Range("E32").Value = '  -2.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.306'
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("E34").Value = '  +1.27%  '
$ws.Range("E35").Value = '  -0.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05742'
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.148'
$ws.Range("E37").Value = '  +1.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02072'
$ws.Range("E38").Value = '  -1.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.720'
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5639'
$ws.Range("E40").Value = '  +0.66%  '
$ws.Range("E41").Value = '  +1.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.885'
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07210'
$ws.Range("E43").Value = '  -1.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.77'
$ws.Range("E44").Value = '  +1.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5283'
$ws.Range("E45").Value = '  +0.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.137'
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.115'
$ws.Range("E47").Value = '  -1.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '113.58'
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.826'
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.409'
$ws.Range("E50").Value = '  +4.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.005'
$ws.Range("E51").Value = '  +0.44%  '
